$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2:A20").Value = "2026-02-18 02:33:59"
